# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet right before the existing "2022-Q1"
#    sheet, populated with a single fund row.
# 2. Insert a new summary row for "2022-Q4" at the top of the "总计" sheet's
#    data (pushing the existing 2022-Q1 / 2021-Q4 rows down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q4" worksheet, positioned before "2022-Q1"
# ---------------------------------------------------------------------
# Duplicate the "2022-Q1" sheet (this keeps its layout/number formats/
# sheet properties intact) and drop it immediately before itself, then
# rename the duplicate. This is more faithful than Worksheets.Add(),
# which creates a brand new blank sheet without the original sheetPr /
# margins / styles.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)

$q4 = $wb.Worksheets.Item("2022-Q1 (2)")
$q4.Name = "2022-Q4"

# The duplicate still carries all 4 fund rows from "2022-Q1" (rows 2-5).
# 2022-Q4 only has a single fund, so drop the extra rows 3-5, keeping
# just the header (row 1) and one data row (row 2).
$q4.Rows.Item(3).Resize(3).Delete()

# Overwrite the remaining data row with the 2022-Q4 fund figures.
# Numeric-looking values are entered with a leading apostrophe so they
# are stored as text (matching the source data, which keeps these as
# text - e.g. to preserve the leading zero in the fund code).
$q4.Range("B2").Value = "'004321"
$q4.Range("C2").Value = "前海开源沪港深强国产业灵活配置混合"
$q4.Range("D2").Value = "'0.33"
$q4.Range("E2").Value = "'90.19"
$q4.Range("F2").Value = "'3.86"
$q4.Range("G2").Value = "'0.0127"
$q4.Range("H2").Value = 10

# ---------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet, inserting a 2022-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Shift the existing data rows down by one (row 3 -> row 4, row 2 -> row 3)
# using Copy so formatting/styles (e.g. the "A" index column style) move
# along with the values.
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# Write the new 2022-Q4 summary values into row 2.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# Fix up the zero-based index column (A2:A4 = 0,1,2).
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
